# no-op
$p = $ppt.ActivePresentation
